$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.700.24'
$ws.Range('E2').Value = '  -2.38%  '
$ws.Range('D3').Value = '1.558.71'
$ws.Range('E3').Value = '  -0.35%  '
$ws.Range('E4').Value = '  -0.04%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '205.93'
$c.Style = "Normal"
$ws.Range('E5').Value = '  -1.11%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '0.490'
$c.Style = "Normal"
$ws.Range('E6').Value = '  -2.03%  '
$ws.Range('E7').Value = '  -0.05%  '
$c = $ws.Range('D8')
$c.NumberFormat = "@"
$c.Value = '21.99'
$c.Style = "Normal"
$ws.Range('E8').Value = '  +0.55%  '
$ws.Range('E9').Value = '  -0.51%  '
$ws.Range('E10').Value = '  -1.52%  '
$ws.Range('E11').Value = '  -0.56%  '
$ws.Range('D12').Value = '1.777.79'
$ws.Range('E12').Value = '  -0.55%  '
$ws.Range('D13').Value = '1.546.65'
$ws.Range('E13').Value = '  -0.99%  '
$ws.Range('E14').Value = '  -2.14%  '
$ws.Range('E15').Value = '  -0.54%  '
$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '61.70'
$c.Style = "Normal"
$ws.Range('E16').Value = '  -2.61%  '
$ws.Range('D17').Value = '26.728.39'
$ws.Range('E17').Value = '  -2.29%  '
$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '213.74'
$c.Style = "Normal"
$ws.Range('E18').Value = '  +0.87%  '
$ws.Range('E19').Value = '  +1.16%  '
$ws.Range('E20').Value = '  -2.01%  '
$ws.Range('E21').Value = '  +0.05%  '
$ws.Range('E22').Value = '  -0.70%  '
$ws.Range('E23').Value = '  -1.60%  '
$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '1.99'
$c.Style = "Normal"
$ws.Range('E24').Value = '  -0.39%  '
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '152.36'
$c.Style = "Normal"
$ws.Range('E25').Value = '  -0.51%  '
$ws.Range('E26').Value = '  +0.90%  '
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '14.83'
$c.Style = "Normal"
$ws.Range('E27').Value = '  -0.92%  '
$ws.Range('E28').Value = '  +0.01%  '
$ws.Range('E29').Value = '  -0.83%  '
$ws.Range('E30').Value = '  -1.41%  '
$ws.Range('E31').Value = '  -4.06%  '
$ws.Range('E32').Value = '  -1.53%  '
$ws.Range('D33').Value = '1.385.68'
$ws.Range('E33').Value = '  +1.63%  '
$ws.Range('E34').Value = '  -1.31%  '
$ws.Range('E35').Value = '  +0.75%  '
$ws.Range('E36').Value = '  -1.00%  '
$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '0.929'
$c.Style = "Normal"
$ws.Range('E37').Value = '  -4.30%  '
$ws.Range('E38').Value = '  -2.43%  '
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '0.812'
$c.Style = "Normal"
$ws.Range('E39').Value = '  -1.02%  '
$ws.Range('B40').Value = 'ImmutableX'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '0.517'
$c.Style = "Normal"
$ws.Range('E40').Value = '  -2.59%  '
$ws.Range('E41').Value = '  +0.03%  '
$ws.Range('E42').Value = '  +2.04%  '
$ws.Range('E43').Value = '  +2.24%  '
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '2.18'
$c.Style = "Normal"
$ws.Range('E44').Value = '  +1.43%  '
$ws.Range('E45').Value = '  -1.37%  '
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '63.12'
$c.Style = "Normal"
$ws.Range('E46').Value = '  -1.30%  '
$ws.Range('D47').Value = '1.692.62'
$ws.Range('E47').Value = '  -0.48%  '
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '85.44'
$c.Style = "Normal"
$ws.Range('E48').Value = '  -0.06%  '
$ws.Range('D49').Value = '0.0₇0974'
$ws.Range('E49').Value = '  -1.84%  '
$ws.Range('E50').Value = '  -0.16%  '
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '0.0948'
$c.Style = "Normal"
$ws.Range('E51').Value = '  -0.65%  '
